$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.781.40'
$ws.Range("E2").Value = '  +0.13%  '

$ws.Range("D3").Value = '2.556.66'
$ws.Range("E3").Value = '  +1.11%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = '311.07'
$ws.Range("E5").Value = '  -1.88%  '

$ws.Range("D6").Value = '98.38'
$ws.Range("E6").Value = '  +1.03%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("E9").Value = '  -0.37%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.60'
$ws.Range("E10").Value = '  -0.56%  '

$ws.Range("E11").Value = '  +0.02%  '

$ws.Range("E12").Value = '  -1.61%  '

$ws.Range("D13").Value = '2.950.28'
$ws.Range("E13").Value = '  +1.28%  '

$ws.Range("E14").Value = '  -1.97%  '

$ws.Range("E15").Value = '  +5.76%  '

$ws.Range("D16").Value = '2.608.97'
$ws.Range("E16").Value = '  +1.18%  '

$ws.Range("E17").Value = '  -1.40%  '

$ws.Range("D18").Value = '42.793.52'
$ws.Range("E18").Value = '  +0.03%  '

$ws.Range("E19").Value = '  -1.28%  '

$ws.Range("D20").Value = '0.0₃0959'
$ws.Range("E20").Value = '  -0.40%  '

$ws.Range("D21").Value = '12.36'
$ws.Range("E21").Value = '  -3.27%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.60'
$ws.Range("E22").Value = '  -0.14%  '

$ws.Range("D23").Value = '247.49'
$ws.Range("E23").Value = '  -1.68%  '

$ws.Range("E24").Value = '  -1.01%  '

$ws.Range("E25").Value = '  -0.08%  '

$ws.Range("D26").Value = '26.77'
$ws.Range("E26").Value = '  +1.04%  '

$ws.Range("E27").Value = '  -0.05%  '

$ws.Range("E28").Value = '  -0.22%  '

$ws.Range("D29").Value = '40.01'
$ws.Range("E29").Value = '  -2.05%  '

$ws.Range("D30").Value = '10.17'
$ws.Range("E30").Value = '  -2.29%  '

$ws.Range("D31").Value = '158.37'
$ws.Range("E31").Value = '  +0.00%  '

$ws.Range("E32").Value = '  -2.96%  '

$ws.Range("E33").Value = '  +0.96%  '

$ws.Range("E34").Value = '  -1.91%  '

$ws.Range("E35").Value = '  -3.24%  '

$ws.Range("B36").Value = 'WEMIXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").Value = '2.64'
$ws.Range("E36").Value = '  -2.99%  '

$ws.Range("B37").Value = 'Celestia'
$ws.Range("C37").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D37").Value = '18.65'
$ws.Range("E37").Value = '  -1.42%  '

$ws.Range("D38").Value = '2.59'
$ws.Range("E38").Value = '  +11.30%  '

$ws.Range("E39").Value = '  -0.74%  '

$ws.Range("D40").Value = '0.118'
$ws.Range("E40").Value = '  -0.64%  '

$ws.Range("D41").Value = '22.64'
$ws.Range("E41").Value = '  +1.19%  '

$ws.Range("E42").Value = '  +6.35%  '

$ws.Range("E43").Value = '  -0.14%  '

$ws.Range("E44").Value = '  -1.17%  '

$ws.Range("D45").Value = '1.989.01'
$ws.Range("E45").Value = '  -2.06%  '

$ws.Range("D46").Value = '3.19'
$ws.Range("E46").Value = '  -2.39%  '

$ws.Range("E47").Value = '  -1.42%  '

$ws.Range("D48").Value = '2.803.60'
$ws.Range("E48").Value = '  +1.35%  '

$ws.Range("D49").Value = '81.27'
$ws.Range("E49").Value = '  -3.79%  '

$ws.Range("E50").Value = '  +1.51%  '

$ws.Range("D51").Value = '73.39'
$ws.Range("E51").Value = '  -2.30%  '
